# Insert a new row at row 106, shifting existing rows 106:197 down to 107:198.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(106).Insert()

# Populate the newly inserted row 106 with the new data record.
$ws.Range("A106").Value = 10
$ws.Range("B106").Value = "Vega Modelo de Temuco"
$ws.Range("C106").Value = "La Araucanía"
$ws.Range("D106").Value = 44447
$ws.Range("E106").Value = 9
$ws.Range("F106").Value = "Fruta"
$ws.Range("G106").Value = 100108
$ws.Range("H106").Value = "Tropicales y subtropicales"
$ws.Range("I106").Value = 100108002
$ws.Range("J106").Value = "Mango"
$ws.Range("K106").Value = "Sin especificar"
$ws.Range("L106").Value = "Primera"
$ws.Range("M106").Value = 200
$ws.Range("N106").Value = 9000
$ws.Range("O106").Value = 9000
$ws.Range("P106").Value = 9000
$ws.Range("Q106").Value = "$/bandeja 4 kilos"
$ws.Range("R106").Value = "Brasil"
$ws.Range("S106").Value = 2250
$ws.Range("T106").Value = 4
